$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 666663244656.6726
    3  = 215894388530.8438
    4  = 84259407156.79712
    5  = 34490373188.7056
    6  = 29812930074.87981
    7  = 25079116536.25563
    8  = 13075337557.07294
    9  = 9938396539.221432
    10 = 9414554630.551594
    11 = 8233547006.034203
    12 = 7345088469.253864
    13 = 6362933528.338446
    14 = 5941956028.061554
    15 = 5606911520.500404
    16 = 5548924012.004469
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $newValues[$row]
}
